$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.341.37'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.74%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.115.99'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.63%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '621.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.01'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +24.94%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.376'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.26%  '

$ws.Range("E9").Value = '  +0.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.113.64'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.75%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.724'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +21.75%  '

$ws.Range("E12").Value = '  +5.85%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.75%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.59'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.44'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.84%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.134.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.693.48'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.100.01'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +12.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000218'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '436.74'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.81'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.55%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.13%  '

$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '87.82'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.20%  '

$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.95%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.281.38'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.54%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("E30").Value = '  -1.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.13'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +12.38%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '525.97'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.901'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -16.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.76'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.02%  '

$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.147'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +13.34%  '

$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.10'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.41%  '

$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '23.72'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.42%  '

$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.30'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.18%  '

$ws.Range("E39").Value = '  +3.43%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0891'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +27.94%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.27'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.08%  '

$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("E43").Value = '  +16.95%  '

$ws.Range("E44").Value = '  +10.04%  '

$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.94'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '149.37'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '44.10'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.31'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.85%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '168.28'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.23'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.13%  '
